$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = 2616
$ws.Range("B2").Value = "MAI"

# --- Row 3 ---
$ws.Range("A3").Value = 1476
$ws.Range("B3").Value = "MAI"

# --- Row 4 ---
$ws.Range("A4").Value = 4781
$ws.Range("B4").Value = "MAI"

# --- Row 5 ---
$ws.Range("A5").Value = 5052
$ws.Range("B5").Value = "MAI"

# G5 must become the new shared string "NO TIME" BEFORE H4 becomes
# "STARTED TOO SOON; MAY BE BAD CURVE" so the shared-string table order
# matches the source workbook (NO TIME=20, STARTED TOO SOON...=21).
$ws.Range("G5").Value = "NO TIME"

# G2-G4 hold time-of-day values, formatted as h:mm (creates a new
# number-format style reusing the existing bordered/wrapped cell style).
$ws.Range("G2").Value = 0.46527777777777773
$ws.Range("G2").NumberFormat = "h:mm"

$ws.Range("G3").Value = 0.52500000000000002
$ws.Range("G3").NumberFormat = "h:mm"

$ws.Range("G4").Value = 0.53749999999999998
$ws.Range("G4").NumberFormat = "h:mm"

$ws.Range("H4").Value = "STARTED TOO SOON; MAY BE BAD CURVE"

# --- Update selection (was A2:B4 / B4, becomes C5) ---
$null = $ws.Range("C5").Select()

# --- Header / footer text updates ---
$ps = $ws.PageSetup
$ps.LeftHeader = "&`"Calibri (Body),Regular`"&24`nDate: JUNE 15, 2023"
$ps.CenterHeader = "&`"Calibri (Body),Regular`"&18Machine (circle):`n&24Ozz   Gib    Alb    Stan    &`"Calibri (Body),Bold`"Yat"
$ps.RightHeader = "&`"Calibri (Body),Regular`"&20`nLicor Enthusiast: EVAN PERKOWSKI"

Write-Host "Edits applied"
